$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 473 (shifts rows 473:506 down to 474:507)
$ws.Rows("473:473").Insert()

# Copy formatting/content of the row above (472) into the newly inserted row 473
$ws.Rows("472:472").Copy()
$ws.Rows("473:473").PasteSpecial()

# Set the new row's specific values
$ws.Range("D473").Value = 45223
$ws.Range("J473").Value = 7000
$ws.Range("K473").Value = 3000
$ws.Range("L473").Value = 3000
$ws.Range("M473").Value = 3000
$ws.Range("P473").Value = 30
